# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text would otherwise be auto-parsed as a number by Excel
# (pure decimals like "1.001" or "0.08930") need an explicit Text format first so
# the literal string - including trailing zeros - is preserved exactly.
$textPriceRows = @(5,6,7,8,9,10,11,12,14,15,16,17,18,19,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51)
foreach ($r in $textPriceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = "27.008.78"
$ws.Cells.Item(2, 5).Value = "  +4.90%  "
$ws.Cells.Item(3, 4).Value = "1.878.60"
$ws.Cells.Item(3, 5).Value = "  +3.74%  "
$ws.Cells.Item(4, 5).Value = "  +0.12%  "
$ws.Cells.Item(5, 4).Value = "279.78"
$ws.Cells.Item(5, 5).Value = "  +1.06%  "
$ws.Cells.Item(6, 4).Value = "1.001"
$ws.Cells.Item(6, 5).Value = "  +0.14%  "
$ws.Cells.Item(7, 4).Value = "0.5281"
$ws.Cells.Item(7, 5).Value = "  +4.34%  "
$ws.Cells.Item(8, 4).Value = "0.3475"
$ws.Cells.Item(8, 5).Value = "  -0.75%  "
$ws.Cells.Item(9, 4).Value = "0.06967"
$ws.Cells.Item(9, 5).Value = "  +4.58%  "
$ws.Cells.Item(10, 4).Value = "20.21"
$ws.Cells.Item(10, 5).Value = "  +1.57%  "
$ws.Cells.Item(11, 4).Value = "0.8116"
$ws.Cells.Item(11, 5).Value = "  -2.05%  "
$ws.Cells.Item(12, 4).Value = "0.07872"
$ws.Cells.Item(12, 5).Value = "  +0.52%  "
$ws.Cells.Item(13, 4).Value = "1.858.51"
$ws.Cells.Item(13, 5).Value = "  +2.60%  "
$ws.Cells.Item(14, 4).Value = "90.05"
$ws.Cells.Item(14, 5).Value = "  +2.84%  "
$ws.Cells.Item(15, 4).Value = "5.159"
$ws.Cells.Item(15, 5).Value = "  +2.12%  "
$ws.Cells.Item(16, 4).Value = "14.54"
$ws.Cells.Item(16, 5).Value = "  +4.32%  "
$ws.Cells.Item(17, 4).Value = "1.000"
$ws.Cells.Item(17, 5).Value = "  +0.07%  "
$ws.Cells.Item(18, 4).Value = "0.000008100"
$ws.Cells.Item(18, 5).Value = "  +0.65%  "
$ws.Cells.Item(19, 4).Value = "1.001"
$ws.Cells.Item(19, 5).Value = "  +0.12%  "
$ws.Cells.Item(20, 4).Value = "27.041.63"
$ws.Cells.Item(20, 5).Value = "  +4.80%  "
$ws.Cells.Item(21, 4).Value = "2.114.68"
$ws.Cells.Item(21, 5).Value = "  +3.84%  "
$ws.Cells.Item(22, 4).Value = "4.753"
$ws.Cells.Item(22, 5).Value = "  +0.87%  "
$ws.Cells.Item(23, 4).Value = "10.07"
$ws.Cells.Item(23, 5).Value = "  +1.05%  "
$ws.Cells.Item(24, 4).Value = "6.182"
$ws.Cells.Item(24, 5).Value = "  +2.40%  "
$ws.Cells.Item(25, 4).Value = "2.352"
$ws.Cells.Item(25, 5).Value = "  +6.59%  "
$ws.Cells.Item(26, 4).Value = "146.62"
$ws.Cells.Item(26, 5).Value = "  +4.98%  "
$ws.Cells.Item(27, 4).Value = "17.47"
$ws.Cells.Item(27, 5).Value = "  +2.78%  "
$ws.Cells.Item(28, 4).Value = "1.667"
$ws.Cells.Item(28, 5).Value = "  +0.13%  "
$ws.Cells.Item(29, 4).Value = "114.57"
$ws.Cells.Item(29, 5).Value = "  +4.80%  "
$ws.Cells.Item(30, 4).Value = "4.363"
$ws.Cells.Item(30, 5).Value = "  +0.91%  "
$ws.Cells.Item(31, 4).Value = "4.351"
$ws.Cells.Item(31, 5).Value = "  +3.13%  "
$ws.Cells.Item(32, 4).Value = "0.08930"
$ws.Cells.Item(32, 5).Value = "  +1.57%  "
$ws.Cells.Item(33, 4).Value = "0.04940"
$ws.Cells.Item(33, 5).Value = "  +1.65%  "
$ws.Cells.Item(34, 4).Value = "1.178"
$ws.Cells.Item(34, 5).Value = "  +3.62%  "
$ws.Cells.Item(35, 4).Value = "0.7350"
$ws.Cells.Item(35, 5).Value = "  +2.08%  "
$ws.Cells.Item(36, 4).Value = "2.895"
$ws.Cells.Item(36, 5).Value = "  +0.54%  "
$ws.Cells.Item(37, 4).Value = "3.311"
$ws.Cells.Item(37, 5).Value = "  +6.67%  "
$ws.Cells.Item(38, 4).Value = "2.386"
$ws.Cells.Item(38, 5).Value = "  +6.82%  "
$ws.Cells.Item(39, 4).Value = "0.01858"
$ws.Cells.Item(39, 5).Value = "  +1.36%  "
$ws.Cells.Item(40, 4).Value = "0.5191"
$ws.Cells.Item(40, 5).Value = "  -0.25%  "
$ws.Cells.Item(41, 4).Value = "0.9641"
$ws.Cells.Item(41, 5).Value = "  +1.17%  "
$ws.Cells.Item(42, 4).Value = "115.48"
$ws.Cells.Item(42, 5).Value = "  +1.95%  "
$ws.Cells.Item(43, 4).Value = "6.188"
$ws.Cells.Item(43, 5).Value = "  +0.83%  "
$ws.Cells.Item(44, 4).Value = "8.102"
$ws.Cells.Item(44, 5).Value = "  +0.85%  "
$ws.Cells.Item(45, 4).Value = "1.001"
$ws.Cells.Item(45, 5).Value = "  +0.16%  "
$ws.Cells.Item(46, 4).Value = "0.4537"
$ws.Cells.Item(46, 5).Value = "  -0.03%  "
$ws.Cells.Item(47, 4).Value = "0.1350"
$ws.Cells.Item(47, 5).Value = "  -0.81%  "
$ws.Cells.Item(48, 4).Value = "9.414"
$ws.Cells.Item(48, 5).Value = "  +0.99%  "
$ws.Cells.Item(49, 4).Value = "36.45"
$ws.Cells.Item(49, 5).Value = "  +0.71%  "
$ws.Cells.Item(50, 4).Value = "1.514"
$ws.Cells.Item(50, 5).Value = "  +0.95%  "
$ws.Cells.Item(51, 4).Value = "0.05934"
$ws.Cells.Item(51, 5).Value = "  +1.92%  "
